$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.852.80'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.510.12'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '604.97'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +4.34%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '171.24'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.70%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.614'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.95%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.507.37'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.30%  '
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  +4.84%  '
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('E12').Value = '  -2.62%  '
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000278'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.082.90'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '619.53'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -8.32%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '8.37'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -4.21%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.513.94'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.829.63'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.18%  '
$ws.Range('E20').Value = '  -2.15%  '
$ws.Range('E21').Value = '  -1.35%  '
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.883'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.94'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -11.29%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '15.74'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.00%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '95.88'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.21%  '
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.59'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.64%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.20'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.68%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '33.14'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.85%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.39'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -4.07%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.06'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.56%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.34'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.35%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.94'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -5.06%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '567.53'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.37%  '
$ws.Range('E36').Value = '  -1.35%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.45'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -4.11%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '57.06'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('E39').Value = '  -3.69%  '
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('E41').Value = '  +2.72%  '
$ws.Range('E42').Value = '  +1.58%  '
$ws.Range('E43').Value = '  -3.64%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.331.07'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.95%  '
$ws.Range('E45').Value = '  +2.75%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '32.99'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.31%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0₃0705'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.61'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('E49').Value = '  -3.77%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '135.20'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.70'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.43%  '

Write-Host "Updated cryptos list"
